# CRUD and Import Employee as User
# Rebuild the employee table: drop Username/Password columns, change header
# "Password" -> "Pangkat", replace the 3-row kepala/Produksi roster with a
# 2-row Neraca/anggota roster, and repoint the Email hyperlinks.

$wb = $excel.ActiveWorkbook

# Start from a clean sheet (same name, same position) so leftover column
# widths / hyperlinks / styles from the old 8-column layout don't linger.
$oldName = $wb.Worksheets.Item(1).Name
$new = $wb.Worksheets.Add()
$new.Name = "__tmp_rebuild__"
$wb.Worksheets.Item($oldName).Delete()
$new.Name = $oldName
$ws = $new

# Header row
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Nama Pegawai"
$ws.Range("C1").Value = "Divisi"
$ws.Range("D1").Value = "NIP"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Pangkat"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Aku"
$ws.Range("C2").Value = "Neraca"
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = "aku@bps.go.id"
$ws.Range("F2").Value = "anggota"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Kamu"
$ws.Range("C3").Value = "Neraca"
$ws.Range("D3").Value = 789
$ws.Range("E3").Value = "kam@bps.go.id"
$ws.Range("F3").Value = "anggota"

# Hyperlinks for the email column
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:aku@bps.go.id")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:kam@bps.go.id")
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"

$ws.Range("G14").Select()
